$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a missing space: "...62.283-A,com o mesmo..." -> "...62.283-A, com o mesmo..."
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("62.283-A,com o me", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "62.283-A, com o me", 2)
Write-Output "Step1 (space before 'com o me'): $found"

# ---------------------------------------------------------------------------
# 2) Remove the duplicated "CLÁUSULA 1ª." text
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("CLÁUSULA 1ª. CLÁUSULA 1ª. ", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "CLÁUSULA 1ª. ", 2)
Write-Output "Step2 (remove duplicate CLAUSULA 1a): $found"

# ---------------------------------------------------------------------------
# 3) Fix typo: "referentes a outra demanda" -> "referente a outra demanda"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("seja referentes a outra demanda", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "seja referente a outra demanda", 2)
Write-Output "Step3 (referentes -> referente): $found"

# ---------------------------------------------------------------------------
# 4) Remove trailing period + stray "b" run:
#    "...entre as partes.b" -> "...entre as partes"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("entre as partes.b", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "entre as partes", 2)
Write-Output "Step4 (remove trailing '.b'): $found"

# ---------------------------------------------------------------------------
# 5) Remove the stray ". " run after "...pela CONTRATADA.. "
#    (original has two runs "...CONTRATADA." + ". " -> should become just
#    "...CONTRATADA." with no trailing space)
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("pela CONTRATADA.. ", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "pela CONTRATADA.", 2)
Write-Output "Step5 (remove stray '. ' run): $found"

# ---------------------------------------------------------------------------
# 6) Remove the _GoBack bookmark
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
    Write-Output "Step6 (deleted _GoBack bookmark): True"
} else {
    Write-Output "Step6 (deleted _GoBack bookmark): not found"
}

Write-Output "Done"
